# Daily attendance processing - 2025-10-21 12:38:09
#
# Normalizes the "Recorded By" column (G): the literal editor "System"
# (exact case) is promoted to the front of the comma-separated editor
# list, the remaining editors keeping their original relative order.
# When no exact "System" entry is present, the editor list is instead
# sorted alphabetically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $raw = $cell.Value2

    if ($null -eq $raw -or $raw -eq "") {
        continue
    }

    $parts = $raw -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -le 1) {
        continue
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.CompareTo("System") -eq 0) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $parts) {
            if ($p.CompareTo("System") -ne 0) {
                $rest += $p
            }
        }
        $newParts = @("System") + $rest
    } else {
        $newParts = $parts | Sort-Object
    }

    $newValue = [string]::Join(", ", $newParts)

    if ($newValue -ne $raw) {
        $cell.Value = $newValue
    }
}
